{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = \"constellation Pers\u00e9e\";\nconst replacement = \"Constellation des G\u00e9meaux\";\nlet newText = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(target) !== -1) {\n    targetPara = para;\n    newText = para.text.split(target).join(replacement);\n    break;\n  }\n}\n\nif (targetPara) {\n  targetPara.clear();\n  targetPara.insertText(newText, Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$target = \"constellation Pers\u00e9e\"\n$replacement = \"Constellation des G\u00e9meaux\"\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $r = $p.Range\n    $full = $r.Text\n    if ($full.Contains($target)) {\n        # Exclude the trailing paragraph mark from the range we rewrite.\n        $r.MoveEnd(1, -1) | Out-Null\n        $newText = $r.Text.Replace($target, $replacement)\n        $r.Delete()\n        $r.InsertAfter($newText)\n        break\n    }\n}\n"}
